# Add I0 and IF columns to the sheet, mirroring the existing IP (column H)
# header style and values: I gets a constant 1, J duplicates H's value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they get the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Headers
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Find the last used row (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, "H").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
